$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.451.18"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").Value = "'3.627.87"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "'591.36"
$ws.Range("E5").Value = "  +1.49%  "
$ws.Range("D6").Value = "'185.40"
$ws.Range("E6").Value = "  +5.57%  "
$ws.Range("D7").Value = "'0.614"
$ws.Range("E7").Value = "  -1.17%  "
$ws.Range("D8").Value = "'0.997"
$ws.Range("E8").Value = "  -0.27%  "
$ws.Range("D9").Value = "'0.675"
$ws.Range("E9").Value = "  -3.58%  "
$ws.Range("E10").Value = "  -8.02%  "
$ws.Range("D11").Value = "'54.29"
$ws.Range("E11").Value = "  -1.01%  "
$ws.Range("D12").Value = "'0.0000255"
$ws.Range("E12").Value = "  -10.50%  "
$ws.Range("D13").Value = "'9.97"
$ws.Range("E13").Value = "  -4.42%  "
$ws.Range("D14").Value = "'4.194.28"
$ws.Range("E14").Value = "  -0.73%  "
$ws.Range("D15").Value = "'3.621.84"
$ws.Range("E15").Value = "  -0.85%  "
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").Value = "'18.46"
$ws.Range("E17").Value = "  -2.90%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "'67.283.06"
$ws.Range("E18").Value = "  -0.39%  "
$ws.Range("D19").Value = "'12.28"
$ws.Range("E19").Value = "  -2.14%  "
$ws.Range("E20").Value = "  -3.21%  "
$ws.Range("D21").Value = "'394.15"
$ws.Range("E21").Value = "  -2.15%  "
$ws.Range("E22").Value = "  -3.09%  "
$ws.Range("D23").Value = "'85.13"
$ws.Range("E23").Value = "  -2.54%  "
$ws.Range("D24").Value = "'2.89"
$ws.Range("E24").Value = "  -2.81%  "
$ws.Range("D25").Value = "'12.32"
$ws.Range("E25").Value = "  -1.54%  "
$ws.Range("D26").Value = "'6.07"
$ws.Range("E26").Value = "  +0.80%  "
$ws.Range("D27").Value = "'10.38"
$ws.Range("E27").Value = "  -1.74%  "
$ws.Range("D28").Value = "'3.62"
$ws.Range("E28").Value = "  -10.44%  "
$ws.Range("D29").Value = "'9.03"
$ws.Range("E29").Value = "  -2.83%  "
$ws.Range("D30").Value = "'31.29"
$ws.Range("E30").Value = "  -2.55%  "
$ws.Range("D31").Value = "'6.82"
$ws.Range("E31").Value = "  -3.61%  "
$ws.Range("B32").Value = "OKB"
$ws.Range("C32").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D32").Value = "'65.56"
$ws.Range("E32").Value = "  +2.47%  "
$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").Value = "'11.93"
$ws.Range("E33").Value = "  -1.81%  "
$ws.Range("D34").Value = "'597.27"
$ws.Range("E34").Value = "  +1.00%  "
$ws.Range("E35").Value = "  -1.91%  "
$ws.Range("D36").Value = "'41.76"
$ws.Range("E36").Value = "  -1.38%  "
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("E38").Value = "  -0.38%  "
$ws.Range("E39").Value = "  -3.86%  "
$ws.Range("D40").Value = "'0.0₃0746"
$ws.Range("E40").Value = "  -14.34%  "
$ws.Range("E41").Value = "  -1.03%  "
$ws.Range("E42").Value = "  -6.19%  "
$ws.Range("E43").Value = "  -3.97%  "
$ws.Range("E44").Value = "  -8.72%  "
$ws.Range("D45").Value = "'2.705.38"
$ws.Range("E45").Value = "  +1.18%  "
$ws.Range("E46").Value = "  -1.89%  "
$ws.Range("D47").Value = "'3.03"
$ws.Range("E47").Value = "  -2.41%  "
$ws.Range("E48").Value = "  -5.66%  "
$ws.Range("D49").Value = "'136.57"
$ws.Range("E49").Value = "  -2.47%  "
$ws.Range("D50").Value = "'8.29"
$ws.Range("E50").Value = "  -6.56%  "
$ws.Range("D51").Value = "'2.59"
$ws.Range("E51").Value = "  -4.43%  "
